$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.898.91"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "3.676.66"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "239.74"
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E6").Value = "  +11.49%  "
$ws.Range("D7").Value = "657.72"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("E9").Value = "  +4.30%  "
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "3.673.13"
$ws.Range("E11").Value = "  +2.40%  "
$ws.Range("D12").Value = "45.61"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "6.83"
$ws.Range("E14").Value = "  +6.55%  "
$ws.Range("D15").Value = "4.364.10"
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("D16").Value = "0.0000270"
$ws.Range("E16").Value = "  +3.93%  "
$ws.Range("D17").Value = "96.602.74"
$ws.Range("E17").Value = "  -0.32%  "
$ws.Range("E18").Value = "  +3.68%  "
$ws.Range("D19").Value = "3.687.31"
$ws.Range("E19").Value = "  +2.84%  "
$ws.Range("D20").Value = "18.74"
$ws.Range("E20").Value = "  +3.94%  "
$ws.Range("D21").Value = "12.77"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("D22").Value = "0.536"
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").Value = "533.56"
$ws.Range("E23").Value = "  +3.69%  "
$ws.Range("D24").Value = "3.52"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "7.20"
$ws.Range("E25").Value = "  +4.50%  "
$ws.Range("D26").Value = "0.0000205"
$ws.Range("E26").Value = "  +0.58%  "
$ws.Range("D27").Value = "102.76"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").Value = "13.55"
$ws.Range("E28").Value = "  +3.86%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.874.08"
$ws.Range("E29").Value = "  +2.40%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.168"
$ws.Range("E30").Value = "  +5.25%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "12.42"
$ws.Range("E31").Value = "  +4.74%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "3.04"
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.90"
$ws.Range("E34").Value = "  +15.07%  "
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").Value = "0.186"
$ws.Range("E35").Value = "  +1.68%  "
$ws.Range("D36").Value = "32.72"
$ws.Range("E36").Value = "  +3.34%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "661.59"
$ws.Range("E37").Value = "  +5.97%  "
$ws.Range("B38").Value = "Binance-PegBSC-USD"
$ws.Range("C38").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "0.597"
$ws.Range("E39").Value = "  +5.41%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "8.90"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.161"
$ws.Range("E41").Value = "  +5.05%  "
$ws.Range("D42").Value = "6.69"
$ws.Range("E42").Value = "  +11.27%  "
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").Value = "2.00"
$ws.Range("E43").Value = "  +2.35%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "0.964"
$ws.Range("E44").Value = "  +4.49%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "38.78"
$ws.Range("E45").Value = "  +17.59%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0463"
$ws.Range("E47").Value = "  +4.70%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.436"
$ws.Range("E48").Value = "  +6.74%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "2.34"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("B50").Value = "MantraDAO"
$ws.Range("C50").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D50").Value = "3.70"
$ws.Range("E50").Value = "  +6.04%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "23.63"
$ws.Range("E51").Value = "  +0.13%  "
